$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new "early_res" column (BI) after the existing "slow_mvmnt" column (BH),
# mirroring the existing "late_res" / "slow_mvmnt" header+description pattern.
$ws.Range("BI1").Value = "early_res"
$ws.Range("BI2").Value = "subject started moving too early."

# Match formatting (wrap text style etc.) used by the rest of the header/description cells.
$ws.Range("BH1:BH2").Copy()
$ws.Range("BI1:BI2").PasteSpecial(-4122)  # xlPasteFormats

# Update the active view to reflect the new last column, like Excel would
# after a user selects/edits the new column.
$ws.Range("BI3").Select()
